# Auto-generated: update 2024 (column K) and, where applicable, 2023 (column J)
# crime-count cells to reflect newly added data for 2024-04-28, across the
# Citywide Totals sheet, the By Neighborhood roll-up sheet, and every affected
# individual neighborhood sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 2332
$ws.Range('K3').Value = 2235
$ws.Range('K4').Value = 470
$ws.Range('K6').Value = 2814
$ws.Range('K7').Value = 7998

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K5').Value = 15
$ws.Range('K6').Value = 63
$ws.Range('K7').Value = 235
$ws.Range('K8').Value = 537
$ws.Range('K9').Value = 32
$ws.Range('K10').Value = 45
$ws.Range('K11').Value = 172
$ws.Range('K14').Value = 46
$ws.Range('J15').Value = 359
$ws.Range('K15').Value = 82
$ws.Range('K18').Value = 53
$ws.Range('K19').Value = 231
$ws.Range('K20').Value = 175
$ws.Range('K25').Value = 34
$ws.Range('K29').Value = 405
$ws.Range('K31').Value = 91
$ws.Range('K33').Value = 314
$ws.Range('K37').Value = 258
$ws.Range('K41').Value = 71
$ws.Range('K42').Value = 276
$ws.Range('K43').Value = 73
$ws.Range('K44').Value = 76
$ws.Range('K49').Value = 54
$ws.Range('K50').Value = 51
$ws.Range('K52').Value = 218
$ws.Range('K54').Value = 152
$ws.Range('K60').Value = 56
$ws.Range('J63').Value = 97
$ws.Range('K64').Value = 53
$ws.Range('K65').Value = 192
$ws.Range('K67').Value = 305
$ws.Range('K73').Value = 79
$ws.Range('K77').Value = 56
$ws.Range('K78').Value = 107
$ws.Range('K79').Value = 210
$ws.Range('K83').Value = 174
$ws.Range('K84').Value = 56
$ws.Range('K85').Value = 391
$ws.Range('K90').Value = 69
$ws.Range('K91').Value = 74
$ws.Range('K95').Value = 126
$ws.Range('K96').Value = 109
$ws.Range('K97').Value = 69
$ws.Range('K99').Value = 144
$ws.Range('K101').Value = 7998

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K2').Value = 18
$ws.Range('K6').Value = 17
$ws.Range('K7').Value = 46

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K2').Value = 37
$ws.Range('K7').Value = 109

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 77
$ws.Range('K6').Value = 62
$ws.Range('K7').Value = 235

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 53
$ws.Range('K3').Value = 43
$ws.Range('K7').Value = 172

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 142
$ws.Range('K3').Value = 131
$ws.Range('K6').Value = 91
$ws.Range('K7').Value = 391

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 59
$ws.Range('K6').Value = 93
$ws.Range('K7').Value = 218

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 159
$ws.Range('K3').Value = 157
$ws.Range('K4').Value = 29
$ws.Range('K7').Value = 537

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 69
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 174

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 117
$ws.Range('K6').Value = 85
$ws.Range('K7').Value = 314

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 126

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 89
$ws.Range('K6').Value = 82
$ws.Range('K7').Value = 258

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K3').Value = 48
$ws.Range('K6').Value = 80
$ws.Range('K7').Value = 192

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K3').Value = 50
$ws.Range('K7').Value = 144

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K3').Value = 16
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 91

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 93
$ws.Range('K7').Value = 305

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K4').Value = 4
$ws.Range('K7').Value = 56

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K2').Value = 6
$ws.Range('K7').Value = 54

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K6').Value = 66
$ws.Range('K7').Value = 152

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 108
$ws.Range('K3').Value = 133
$ws.Range('K6').Value = 133
$ws.Range('K7').Value = 405

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K3').Value = 65
$ws.Range('K6').Value = 74
$ws.Range('K7').Value = 231

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K6').Value = 37
$ws.Range('K7').Value = 76

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K3').Value = 21
$ws.Range('K7').Value = 63

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('K2').Value = 25
$ws.Range('K7').Value = 71

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 69
$ws.Range('K3').Value = 82
$ws.Range('K6').Value = 113
$ws.Range('K7').Value = 276

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K6').Value = 23
$ws.Range('K7').Value = 45

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K3').Value = 27
$ws.Range('K4').Value = 8
$ws.Range('K7').Value = 107

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 31
$ws.Range('K7').Value = 74

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K6').Value = 50
$ws.Range('K7').Value = 210

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K6').Value = 17
$ws.Range('K7').Value = 53

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 58
$ws.Range('K7').Value = 175

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K2').Value = 18
$ws.Range('K7').Value = 53

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K3').Value = 13
$ws.Range('K7').Value = 34

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K2').Value = 27
$ws.Range('K3').Value = 20
$ws.Range('J5').Value = 8
$ws.Range('K6').Value = 29
$ws.Range('J7').Value = 359
$ws.Range('K7').Value = 82

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('K4').Value = 5
$ws.Range('K7').Value = 51

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('K2').Value = 10
$ws.Range('K7').Value = 32

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 79

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K2').Value = 14
$ws.Range('K7').Value = 69

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('K2').Value = 6
$ws.Range('K7').Value = 15

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K2').Value = 27
$ws.Range('K6').Value = 18
$ws.Range('K7').Value = 69

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K2').Value = 15
$ws.Range('K3').Value = 21
$ws.Range('K7').Value = 56

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('K3').Value = 20
$ws.Range('K7').Value = 73

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K2').Value = 26
$ws.Range('K7').Value = 56
